$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the sentence that needs to be split out into its own bullet:
#   "With Memori, I intend ... when travelling. As there is no current
#    renowned alternative."
# We search for just the trailing sentence so the match is specific
# and unambiguous.
# ------------------------------------------------------------------
$oldSentence = ". As there is no current renowned alternative."

$match = $d.Content.Duplicate
$found = $match.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Remember which paragraph (by 1-based index) holds the match, so we
    # can reliably reach the freshly-created paragraph afterwards too.
    $beforeMatch = $d.Range(0, $match.Start)
    $paraIndex = $beforeMatch.Paragraphs.Count

    # Split the paragraph in two right after the sentence we found -
    # this turns the tail ("As it stands ... first of its kind") into
    # its own list item, inheriting the same ListParagraph/numbering
    # formatting as the paragraph it was split from.
    $splitPoint = $d.Range($match.End, $match.End)
    $splitPoint.InsertParagraphAfter()

    # Collapse the old sentence down to a single trailing space, so the
    # first paragraph now just ends with "... when travelling ".
    $oldRange = $d.Range($match.Start, $match.End)
    $oldRange.Text = " "

    # Populate the newly created paragraph with the new wording.
    $newPara = $d.Paragraphs($paraIndex + 1)
    $ip = $d.Range($newPara.Range.Start, $newPara.Range.Start)

    $ip.InsertAfter("As")
    $ip.Collapse(0)
    $ip.InsertAfter(" it stands there is no current")
    $ip.Collapse(0)
    $ip.InsertAfter(" ")
    $ip.Collapse(0)
    $ip.InsertAfter("renowned alternative")
    $ip.Collapse(0)
    $ip.InsertAfter(" for this, hence why I believe this will be the first of its kind")
}
